# "unify the conception of DataNode, DataTable, Entity."
# The sheet previously named "Property1" is renamed to "DataNode" to align
# with the unified DataNode/DataTable/Entity terminology used elsewhere in
# the project.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")
$ws.Name = "DataNode"

# Move the active cell/selection to where the author last left it (W37) on
# the renamed sheet.
[void]$ws.Range("W37").Select()
